# [BI-1613] Update TAF to include term type
# Adds a new "Term Type" column (R) to the Template sheet, with the same
# header formatting (bold, wrap text, bordered) used by the rest of the
# header row, and moves the selection onto the new header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$newHeaderCell = $ws.Range("R1")
$newHeaderCell.Value = "Term Type"

# Match formatting used for the rest of row 1's headers: bold 11pt font,
# wrapped text, and the thin themed border used throughout the header row.
$newHeaderCell.Font.Bold = $true
$newHeaderCell.Font.Size = 11
$newHeaderCell.WrapText = $true
[void]($newHeaderCell.Borders.LineStyle = 1)

[void]($newHeaderCell.Select())
